# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.386.74"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.883.36"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'0.7121"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'242.37"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07987"
$ws.Range("E8").Value = "  +2.86%  "
$ws.Range("D9").Value = "'0.3126"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "'25.24"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "'0.08338"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "1.893.46"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "'0.7203"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "'5.243"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'92.60"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "'6.305"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("E17").Value = "  +2.36%  "
$ws.Range("D18").Value = "29.398.37"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'240.95"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "2.139.94"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'7.838"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'0.1587"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("D27").Value = "'9.072"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").Value = "'18.57"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "'1.507"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'4.416"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "'4.339"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'1.205"
$ws.Range("E32").Value = "  -5.73%  "
$ws.Range("D33").Value = "'0.05371"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "'1.950"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "'1.182"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "'0.7486"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D39").Value = "1.289.96"
$ws.Range("E39").Value = "  +9.89%  "
$ws.Range("D40").Value = "'2.746"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").Value = "'6.601"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").Value = "'0.9103"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value = "'111.60"
$ws.Range("E43").Value = "  +4.86%  "
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("D47").Value = "2.036.62"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("D48").Value = "'1.809"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'0.5219"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'9.492"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'0.4394"
$ws.Range("E51").Value = "  +1.95%  "
